# Weekly update to the "Hortaliza, Femacal de La Calera - Papa" data sheet.
# Three new weekly price records are inserted above the existing row 448
# (pushing the former rows 448-453 down to 451-456, unchanged), and the
# three newly inserted rows (448-450) are populated with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 448; existing rows 448:453 shift to 451:456.
$ws.Rows("448:450").Insert()

function Set-DataRow {
    param($r, $D, $H, $I, $J, $K, $L, $M, $N, $O, $P)

    $ws.Cells.Item($r, 1).Value2 = 3
    $ws.Cells.Item($r, 2).Value2 = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value2 = "Coquimbo"
    $ws.Cells.Item($r, 4).Value2 = $D
    $ws.Cells.Item($r, 5).Value2 = 5
    $ws.Cells.Item($r, 6).Value2 = 100114001
    $ws.Cells.Item($r, 7).Value2 = "Papa"
    $ws.Cells.Item($r, 8).Value2 = $H
    $ws.Cells.Item($r, 9).Value2 = $I
    $ws.Cells.Item($r, 10).Value2 = $J
    $ws.Cells.Item($r, 11).Value2 = $K
    $ws.Cells.Item($r, 12).Value2 = $L
    $ws.Cells.Item($r, 13).Value2 = $M
    $ws.Cells.Item($r, 14).Value2 = $N
    $ws.Cells.Item($r, 15).Value2 = $O
    $ws.Cells.Item($r, 16).Value2 = $P
    $ws.Cells.Item($r, 17).Value2 = 25
    $ws.Cells.Item($r, 18).Value2 = "Hortaliza"
}

Set-DataRow 448 44595 "Asterix" "1a (cosecha)" 160 7500 7500 7500 "`$/saco 25 kilos" "Provincia de Talca"    300
Set-DataRow 449 44595 "Rosara"  "1a (cosecha)" 400 7500 8000 7775 "`$/saco 25 kilos" "Provincia de Quillota" 311
Set-DataRow 450 44595 "Rosara"  "1a (cosecha)" 510 7500 8000 7745 "`$/saco 25 kilos" "Provincia de Talca"    310
